$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INCO")

# Delete the entire row 50 (case -282). This shifts row 51 (case -389) up to
# become the new row 50, matching the data in the diff.
$ws.Rows("50:50").Delete()
